# Auto-generated PowerShell Excel COM-interop script
# Applies the cryptos list update described in the commit diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.256.86"
$ws.Range("E2").Value = "  +0.73%  "
$ws.Range("D3").Value = "3.347.06"
$ws.Range("E3").Value = "  +0.40%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "583.02"
$ws.Range("E5").Value = "  +0.08%  "
$ws.Range("D6").Value = "177.09"
$ws.Range("E6").Value = "  +0.68%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("E8").Value = "  +0.09%  "
$ws.Range("E9").Value = "  +2.73%  "
$ws.Range("D10").Value = "0.581"
$ws.Range("E10").Value = "  +0.78%  "
$ws.Range("D11").Value = "47.91"
$ws.Range("E11").Value = "  +5.07%  "
$ws.Range("E12").Value = "  +1.12%  "
$ws.Range("D13").Value = "684.03"
$ws.Range("E13").Value = "  +3.82%  "
$ws.Range("D14").Value = "3.892.86"
$ws.Range("E14").Value = "  +0.57%  "
$ws.Range("D15").Value = "8.40"
$ws.Range("E15").Value = "  +0.00%  "
$ws.Range("D16").Value = "68.314.55"
$ws.Range("E16").Value = "  +0.61%  "
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "3.377.37"
$ws.Range("E17").Value = "  +1.28%  "
$ws.Range("B18").Value = "TRON"
$ws.Range("C18").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D18").Value = "0.119"
$ws.Range("E18").Value = "  +0.99%  "
$ws.Range("D19").Value = "17.42"
$ws.Range("E19").Value = "  -0.01%  "
$ws.Range("D20").Value = "11.17"
$ws.Range("E20").Value = "  +2.03%  "
$ws.Range("D21").Value = "0.893"
$ws.Range("E21").Value = "  +0.41%  "
$ws.Range("D22").Value = "5.43"
$ws.Range("E22").Value = "  +0.19%  "
$ws.Range("D23").Value = "16.89"
$ws.Range("E23").Value = "  -0.95%  "
$ws.Range("D24").Value = "99.86"
$ws.Range("E24").Value = "  +0.32%  "
$ws.Range("E25").Value = "  +1.30%  "
$ws.Range("D26").Value = "2.69"
$ws.Range("E26").Value = "  +0.68%  "
$ws.Range("D27").Value = "9.51"
$ws.Range("E27").Value = "  +2.60%  "
$ws.Range("E28").Value = "  -1.73%  "
$ws.Range("D29").Value = "8.49"
$ws.Range("E29").Value = "  +0.66%  "
$ws.Range("E30").Value = "  -6.69%  "
$ws.Range("D31").Value = "562.16"
$ws.Range("E31").Value = "  -4.92%  "
$ws.Range("D32").Value = "11.05"
$ws.Range("E32").Value = "  +0.75%  "
$ws.Range("E33").Value = "  +0.59%  "
$ws.Range("D34").Value = "57.87"
$ws.Range("E34").Value = "  +1.92%  "
$ws.Range("E35").Value = "  +0.01%  "
$ws.Range("D36").Value = "3.702.48"
$ws.Range("E36").Value = "  -0.50%  "
$ws.Range("D37").Value = "3.28"
$ws.Range("E37").Value = "  -2.12%  "
$ws.Range("E38").Value = "  +3.87%  "
$ws.Range("D39").Value = "34.63"
$ws.Range("E39").Value = "  +2.11%  "
$ws.Range("D40").Value = "3.16"
$ws.Range("E40").Value = "  +1.24%  "
$ws.Range("D41").Value = "2.60"
$ws.Range("E41").Value = "  -1.01%  "
$ws.Range("D42").Value = "0.335"
$ws.Range("E42").Value = "  +0.43%  "
$ws.Range("D43").Value = "0.0₃0670"
$ws.Range("E43").Value = "  +0.65%  "
$ws.Range("D44").Value = "3.25"
$ws.Range("E44").Value = "  -0.37%  "
$ws.Range("D45").Value = "0.0410"
$ws.Range("E45").Value = "  +0.63%  "
$ws.Range("D46").Value = "2.64"
$ws.Range("E46").Value = "  +2.04%  "
$ws.Range("D47").Value = "0.128"
$ws.Range("E47").Value = "  +0.33%  "
$ws.Range("E48").Value = "  -0.13%  "
$ws.Range("E49").Value = "  -0.67%  "
$ws.Range("D50").Value = "131.24"
$ws.Range("E50").Value = "  +3.13%  "
$ws.Range("E51").Value = "  -0.48%  "
